$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.069.21"
$ws.Range("E2").Value = "  +2.02%  "
$ws.Range("D3").Value = "1.910.53"
$ws.Range("E3").Value = "  +2.19%  "
$ws.Range("D4").Value = "'1.006"
$ws.Range("E4").Value = "  -0.77%  "
$ws.Range("D5").Value = "'316.78"
$ws.Range("E5").Value = "  +1.35%  "
$ws.Range("E6").Value = "  -0.89%  "
$ws.Range("D7").Value = "'0.4811"
$ws.Range("E7").Value = "  +0.69%  "
$ws.Range("D8").Value = "'0.3815"
$ws.Range("E8").Value = "  +0.95%  "
$ws.Range("D9").Value = "'0.07352"
$ws.Range("E9").Value = "  -0.18%  "
$ws.Range("D10").Value = "'0.9348"
$ws.Range("E10").Value = "  -0.22%  "
$ws.Range("D11").Value = "'20.83"
$ws.Range("E11").Value = "  +0.35%  "
$ws.Range("D12").Value = "'0.07771"
$ws.Range("E12").Value = "  -0.92%  "
$ws.Range("D13").Value = "1.875.32"
$ws.Range("E13").Value = "  +0.22%  "
$ws.Range("D14").Value = "'5.504"
$ws.Range("E14").Value = "  +1.22%  "
$ws.Range("D15").Value = "'6.624"
$ws.Range("E15").Value = "  +0.66%  "
$ws.Range("E16").Value = "  +0.95%  "
$ws.Range("E17").Value = "  -0.80%  "
$ws.Range("D18").Value = "'0.000008829"
$ws.Range("E18").Value = "  -0.96%  "
$ws.Range("E19").Value = "  -0.79%  "
$ws.Range("D20").Value = "28.104.61"
$ws.Range("E20").Value = "  +2.05%  "
$ws.Range("D21").Value = "'14.85"
$ws.Range("E21").Value = "  -0.65%  "
$ws.Range("D22").Value = "'5.172"
$ws.Range("E22").Value = "  +0.71%  "
$ws.Range("D23").Value = "2.135.53"
$ws.Range("E23").Value = "  +1.24%  "
$ws.Range("D24").Value = "'10.92"
$ws.Range("E24").Value = "  +1.73%  "
$ws.Range("D25").Value = "'156.02"
$ws.Range("E25").Value = "  +1.32%  "
$ws.Range("D26").Value = "'1.918"
$ws.Range("E26").Value = "  -2.33%  "
$ws.Range("E27").Value = "  -0.07%  "
$ws.Range("D28").Value = "'2.117"
$ws.Range("E28").Value = "  +4.58%  "
$ws.Range("D29").Value = "'116.54"
$ws.Range("E29").Value = "  +0.54%  "
$ws.Range("D30").Value = "'4.956"
$ws.Range("E30").Value = "  -0.89%  "
$ws.Range("D31").Value = "'0.08947"
$ws.Range("E31").Value = "  +0.19%  "
$ws.Range("D32").Value = "'3.319"
$ws.Range("E32").Value = "  -0.57%  "
$ws.Range("D33").Value = "'1.256"
$ws.Range("E33").Value = "  +2.81%  "
$ws.Range("D34").Value = "'0.7753"
$ws.Range("E34").Value = "  +2.66%  "
$ws.Range("D35").Value = "'4.680"
$ws.Range("E35").Value = "  +1.49%  "
$ws.Range("D36").Value = "'2.651"
$ws.Range("E36").Value = "  -1.73%  "
$ws.Range("D37").Value = "'0.02058"
$ws.Range("E37").Value = "  +0.23%  "
$ws.Range("E38").Value = "  -0.74%  "
$ws.Range("D39").Value = "'0.05317"
$ws.Range("E39").Value = "  +0.75%  "
$ws.Range("D40").Value = "'0.5489"
$ws.Range("E40").Value = "  +2.34%  "
$ws.Range("D41").Value = "'2.988"
$ws.Range("E41").Value = "  -0.27%  "
$ws.Range("D42").Value = "'7.018"
$ws.Range("E42").Value = "  -0.94%  "
$ws.Range("D43").Value = "'0.1528"
$ws.Range("E43").Value = "  +0.20%  "
$ws.Range("D44").Value = "'8.490"
$ws.Range("E44").Value = "  +0.36%  "
$ws.Range("D45").Value = "'10.77"
$ws.Range("E45").Value = "  +0.79%  "
$ws.Range("D46").Value = "'0.4835"
$ws.Range("E46").Value = "  +0.35%  "
$ws.Range("E47").Value = "  +5.43%  "
$ws.Range("E48").Value = "  -1.03%  "
$ws.Range("D49").Value = "'1.656"
$ws.Range("E49").Value = "  -0.31%  "
$ws.Range("D50").Value = "'68.12"
$ws.Range("E50").Value = "  +0.84%  "
$ws.Range("D51").Value = "'0.06085"
$ws.Range("E51").Value = "  -0.01%  "
